$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$cellData = @(
    @(1, 1, "16 x 27", "  2    7", "1|    |", "6|    |"),
    @(1, 2, "52 x 84", "  8    4", "5|    |", "2|    |"),
    @(1, 3, "34 x 74", "  7    4", "3|    |", "4|    |"),
    @(2, 1, "76 x 29", "  2    9", "7|    |", "6|    |"),
    @(2, 2, "12 x 30", "  3    0", "1|    |", "2|    |"),
    @(2, 3, "97 x 95", "  9    5", "9|    |", "7|    |"),
    @(3, 1, "75 x 76", "  7    6", "7|    |", "5|    |"),
    @(3, 2, "28 x 92", "  9    2", "2|    |", "8|    |"),
    @(3, 3, "17 x 45", "  4    5", "1|    |", "7|    |"),
    @(4, 1, "16 x 16", "  1    6", "1|    |", "6|    |"),
    @(4, 2, "91 x 53", "  5    3", "9|    |", "1|    |"),
    @(4, 3, "69 x 41", "  4    1", "6|    |", "9|    |"),
    @(5, 1, "29 x 61", "  6    1", "2|    |", "9|    |"),
    @(5, 2, "36 x 13", "  1    3", "3|    |", "6|    |"),
    @(5, 3, "66 x 38", "  3    8", "6|    |", "6|    |")
)

foreach ($entry in $cellData) {
    $row = $entry[0]
    $col = $entry[1]
    $title = $entry[2]
    $digits = $entry[3]
    $corner1 = $entry[4]
    $corner2 = $entry[5]
    $cell = $t.Cell($row, $col)
    $newText = $title + $vtab + $digits + $vtab + "  ----" + $vtab + $corner1 + $vtab + $corner2
    $cell.Range.Text = $newText
}

Write-Host "Updated $($cellData.Count) cells"
